$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)      # 2021-Q3 (stays the originally active sheet)
$ws2 = $wb.Worksheets.Item(2)      # 2021-Q4 (style donor for the fund-holdings layout)
$wsTotal = $wb.Worksheets.Item(3)  # currently "总计" (header + 2 summary rows)

# Rename the existing "总计" sheet to "2022-Q1" FIRST so the name "总计" becomes free
# for the brand-new sheet created below.
$wsTotal.Name = "2022-Q1"

# ---- 1) Create the brand-new "总计" sheet at the end of the workbook, seeded from
#         the content that used to live in the old "总计" sheet. ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNewTotal = $wb.Worksheets.Add($null, $lastSheet)
$wsNewTotal.Name = "总计"

$wsTotal.Range("B1:D1").Copy($wsNewTotal.Range("B1"))

# The existing data rows (2,3) move down to rows 3,4 to make room for the new
# "2022-Q1" entry in row 2.
$wsTotal.Range("B2:D3").Copy($wsNewTotal.Range("B3"))
$wsTotal.Range("A2:A3").Copy($wsNewTotal.Range("A3"))

# Seed row 2's formatting from row 3 (same "A" style, no special style for B:D).
$wsNewTotal.Range("A3").Copy($wsNewTotal.Range("A2"))

$wsNewTotal.Range("A2").Value = 0
$wsNewTotal.Range("B2").Value = "2022-Q1"
$wsNewTotal.Range("C2").Value = 2
$wsNewTotal.Range("D2").Value = 0.32

$wsNewTotal.Range("A3").Value = 1
$wsNewTotal.Range("A4").Value = 2

# ---- 2) Turn the renamed "2022-Q1" sheet into a fund-holdings sheet (same layout
#         used by the 2021-Q3 / 2021-Q4 sheets). ----
$wsTotal.Cells.Clear()

# Copy header formatting/style from the 2021-Q4 sheet (already the fund-holdings layout).
$ws2.Range("B1:H1").Copy($wsTotal.Range("B1"))
$ws2.Range("A2:A3").Copy($wsTotal.Range("A2"))

$wsTotal.Range("B1").Value = "基金代码"
$wsTotal.Range("C1").Value = "基金名称"
$wsTotal.Range("D1").Value = "基金规模"
$wsTotal.Range("E1").Value = "股票总仓位"
$wsTotal.Range("F1").Value = "仓位占比"
$wsTotal.Range("G1").Value = "持有市值(亿元)"
$wsTotal.Range("H1").Value = "仓位排名"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1

# Force the data columns to be stored as plain text (matching the source data, which
# keeps things like leading zeros in fund codes and fixed-looking decimals as text
# rather than numbers). Using a self-quoting formula and then flattening it to a
# static value via PasteSpecial(values) avoids leaving any numeric NumberFormat
# style behind on these cells.
$wsTotal.Range("B2").Formula = '="004738"'
$wsTotal.Range("C2").Formula = '="上投摩根安隆回报混合A"'
$wsTotal.Range("D2").Formula = '="23.04"'
$wsTotal.Range("E2").Formula = '="21.18"'
$wsTotal.Range("F2").Formula = '="1.07"'
$wsTotal.Range("G2").Formula = '="0.2465"'
$wsTotal.Range("H2").Value = 8

$wsTotal.Range("B3").Formula = '="004739"'
$wsTotal.Range("C3").Formula = '="上投摩根安隆回报混合C"'
$wsTotal.Range("D3").Formula = '="7.32"'
$wsTotal.Range("E3").Formula = '="21.18"'
$wsTotal.Range("F3").Formula = '="1.07"'
$wsTotal.Range("G3").Formula = '="0.0783"'
$wsTotal.Range("H3").Value = 8

$txt = $wsTotal.Range("B2:G3")
$txt.Copy()
$txt.PasteSpecial(-4163)

# Restore the originally active sheet/tab.
$ws1.Activate()
